$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.4304097817472439
$ws.Range("C2").Value = 0.09160879547007994
$ws.Range("D2").Value = 0.6619293646000415
$ws.Range("E2").Value = 0.2699344019209917
$ws.Range("G2").Value = 0.002501922895753041
$ws.Range("I2").Value = 0.9731041415514596
$ws.Range("J2").Value = 0.1404529117855091
$ws.Range("K2").Value = 0.5046618618806065
$ws.Range("O2").Value = 4.988518476323748

$ws.Range("B3").Value = 0.3938766318304943
$ws.Range("C3").Value = 0.0841942546234975
$ws.Range("D3").Value = 0.6518000325876869
$ws.Range("E3").Value = 0.2649531102901932
$ws.Range("G3").Value = 0.002504926366483988
$ws.Range("I3").Value = 0.9797275938363192
$ws.Range("J3").Value = 0.1371430818120984
$ws.Range("K3").Value = 0.4623117716852505
$ws.Range("O3").Value = 5.007478707090172

$ws.Range("B4").Value = 0.3715538297486205
$ws.Range("C4").Value = 0.07966748471220342
$ws.Range("D4").Value = 0.6459174939957109
$ws.Range("E4").Value = 0.2620383396620625
$ws.Range("G4").Value = 0.002506868784855129
$ws.Range("I4").Value = 0.9842666384750558
$ws.Range("J4").Value = 0.1351900985987129
$ws.Range("K4").Value = 0.4364395291334517
$ws.Range("O4").Value = 5.021647956366991

$ws.Range("B5").Value = 0.3624849130641792
$ws.Range("C5").Value = 0.0778293269259791
$ws.Range("D5").Value = 0.6436051109912739
$ws.Range("E5").Value = 0.2608867114640745
$ws.Range("G5").Value = 0.002507685123758067
$ws.Range("I5").Value = 0.9862350214723676
$ws.Range("J5").Value = 0.1344141646554107
$ws.Range("K5").Value = 0.4259297850132953
$ws.Range("O5").Value = 5.028057347850762

$ws.Range("B6").Value = 0.360980719800807
$ws.Range("C6").Value = 0.07752449867093958
$ws.Range("D6").Value = 0.6432262656077512
$ws.Range("E6").Value = 0.2606976691631075
$ws.Range("G6").Value = 0.002507822175736469
$ws.Range("I6").Value = 0.9865690366022797
$ws.Range("J6").Value = 0.1342865244243256
$ws.Range("K6").Value = 0.4241866821832332
$ws.Range("O6").Value = 5.029159985438696

$ws.Range("B7").Value = 0.3714314099493379
$ws.Range("C7").Value = 0.07964266813225152
$ws.Range("D7").Value = 0.6458859649338251
$ws.Range("E7").Value = 0.2620226619530186
$ws.Range("G7").Value = 0.002506879693858781
$ws.Range("I7").Value = 0.9842927043028382
$ws.Range("J7").Value = 0.135179553427939
$ws.Range("K7").Value = 0.4362976551376505
$ws.Range("O7").Value = 5.021731823706887

$ws.Range("B8").Value = 0.4177908620065125
$ws.Range("C8").Value = 0.08904693868323932
$ws.Range("D8").Value = 0.6583668876826607
$ws.Range("E8").Value = 0.268187020837523
$ws.Range("G8").Value = 0.002502938141145274
$ws.Range("I8").Value = 0.9752898638479159
$ws.Range("J8").Value = 0.139295223673173
$ws.Range("K8").Value = 0.4900326667066111
$ws.Range("O8").Value = 4.994531188641957

$ws.Range("B9").Value = 0.5095481846159089
$ws.Range("C9").Value = 0.1076919841284507
$ws.Range("D9").Value = 0.6855137895918233
$ws.Range("E9").Value = 0.2814163662761686
$ws.Range("G9").Value = 0.002495985080678752
$ws.Range("I9").Value = 0.9613850768520322
$ws.Range("J9").Value = 0.1479960157967
$ws.Range("K9").Value = 0.5964293871240613
$ws.Range("O9").Value = 4.961262368460268

$ws.Range("B10").Value = 0.5774642032066311
$ws.Range("C10").Value = 0.1215141631748224
$ws.Range("D10").Value = 0.7070887049854377
$ws.Range("E10").Value = 0.2918335436431079
$ws.Range("G10").Value = 0.002491345080763991
$ws.Range("I10").Value = 0.9534596805762874
$ws.Range("J10").Value = 0.154774854308954
$ws.Range("K10").Value = 0.6752093183136196
$ws.Range("O10").Value = 4.949083903274925

$ws.Range("B11").Value = 0.6084674284220171
$ws.Range("C11").Value = 0.1278291295677718
$ws.Range("D11").Value = 0.7172581332750667
$ws.Range("E11").Value = 0.296724576167442
$ws.Range("G11").Value = 0.002489334904307751
$ws.Range("I11").Value = 0.9503525808960589
$ws.Range("J11").Value = 0.1579431740546937
$ws.Range("K11").Value = 0.7111786567213016
$ws.Range("O11").Value = 4.946213117703991

$ws.Range("B12").Value = 0.6202226907077488
$ws.Range("C12").Value = 0.1302243239365168
$ws.Range("D12").Value = 0.7211600387744852
$ws.Range("E12").Value = 0.2985985819411283
$ws.Range("G12").Value = 0.002488588088363201
$ws.Range("I12").Value = 0.9492477253980596
$ws.Range("J12").Value = 0.1591551249500895
$ws.Range("K12").Value = 0.7248179172146649
$ws.Range("O12").Value = 4.945510311106545

$ws.Range("B13").Value = 0.6176903248359338
$ws.Range("C13").Value = 0.1297083054640211
$ws.Range("D13").Value = 0.7203174280816427
$ws.Range("E13").Value = 0.2981940083853161
$ws.Range("G13").Value = 0.002488748289622698
$ws.Range("I13").Value = 0.9494824838185991
$ws.Range("J13").Value = 0.1588935678623358
$ws.Range("K13").Value = 0.7218796460117289
$ws.Range("O13").Value = 4.945644573967542

$ws.Range("B14").Value = 0.6094342428811501
$ws.Range("C14").Value = 0.1280261067344384
$ws.Range("D14").Value = 0.7175781243932704
$ws.Range("E14").Value = 0.2968783133360731
$ws.Range("G14").Value = 0.00248927317523819
$ws.Range("I14").Value = 0.9502602455345794
$ws.Range("J14").Value = 0.1580426378660604
$ws.Range("K14").Value = 0.7123004000476669
$ws.Range("O14").Value = 4.946147592835757

$ws.Range("B15").Value = 0.6043790947363448
$ws.Range("C15").Value = 0.1269962116517434
$ws.Range("D15").Value = 0.7159068567263205
$ws.Range("E15").Value = 0.2960752610097828
$ws.Range("G15").Value = 0.002489596555450855
$ws.Range("I15").Value = 0.9507459919702654
$ws.Range("J15").Value = 0.1575230048543119
$ws.Range("K15").Value = 0.7064352229577651
$ws.Range("O15").Value = 4.94650576703927

$ws.Range("B16").Value = 0.5754401976520853
$ws.Range("C16").Value = 0.1211020075614613
$ws.Range("D16").Value = 0.7064312443167182
$ws.Range("E16").Value = 0.2915169651549832
$ws.Range("G16").Value = 0.002491478468469466
$ws.Range("I16").Value = 0.9536727700525063
$ws.Range("J16").Value = 0.1545695004569723
$ws.Range("K16").Value = 0.6728612479742537
$ws.Range("O16").Value = 4.949325267566934

$ws.Range("B17").Value = 0.5577143971493683
$ws.Range("C17").Value = 0.117493027697094
$ws.Range("D17").Value = 0.7007091096990621
$ws.Range("E17").Value = 0.2887595753682248
$ws.Range("G17").Value = 0.00249265867349302
$ws.Range("I17").Value = 0.9555959225973751
$ws.Range("J17").Value = 0.1527792962569947
$ws.Range("K17").Value = 0.652298125407782
$ws.Range("O17").Value = 4.951738928695448

$ws.Range("B18").Value = 0.5475291701259266
$ws.Range("C18").Value = 0.1154198002774365
$ws.Range("D18").Value = 0.6974512973398532
$ws.Range("E18").Value = 0.2871879246056821
$ws.Range("G18").Value = 0.002493346967579536
$ws.Range("I18").Value = 0.9567489565829632
$ws.Range("J18").Value = 0.1517575795941752
$ws.Range("K18").Value = 0.6404832121396851
$ws.Range("O18").Value = 4.953378418276174

$ws.Range("B19").Value = 0.544082391777863
$ws.Range("C19").Value = 0.114718282897087
$ws.Range("D19").Value = 0.6963539980486075
$ws.Range("E19").Value = 0.2866582516062053
$ws.Range("G19").Value = 0.00249358164088191
$ws.Range("I19").Value = 0.9571474041771424
$ws.Range("J19").Value = 0.1514130108868841
$ws.Range("K19").Value = 0.6364850399465354
$ws.Range("O19").Value = 4.953976653792324

$ws.Range("B20").Value = 0.5596002884411462
$ws.Range("C20").Value = 0.1178769449788319
$ws.Range("D20").Value = 0.7013147835872644
$ws.Range("E20").Value = 0.289051621504683
$ws.Range("G20").Value = 0.002492532058708363
$ws.Range("I20").Value = 0.9553863463566188
$ws.Range("J20").Value = 0.1529690424797536
$ws.Range("K20").Value = 0.6544858200752515
$ws.Range("O20").Value = 4.951455987131851

$ws.Range("B21").Value = 0.6118588507835341
$ws.Range("C21").Value = 0.1285201051334752
$ws.Range("D21").Value = 0.7183813416962721
$ws.Range("E21").Value = 0.2972641711767778
$ws.Range("G21").Value = 0.002489118613440364
$ws.Range("I21").Value = 0.9500298504210889
$ws.Range("J21").Value = 0.1582922460512606
$ws.Range("K21").Value = 0.7151135601321528
$ws.Range("O21").Value = 4.945989410584986

$ws.Range("B22").Value = 0.6461001240931807
$ws.Range("C22").Value = 0.1354984621241613
$ws.Range("D22").Value = 0.72983234788137
$ws.Range("E22").Value = 0.3027590642418616
$ws.Range("G22").Value = 0.002486971599042342
$ws.Range("I22").Value = 0.9469472456289978
$ws.Range("J22").Value = 0.1618422500051366
$ws.Range("K22").Value = 0.7548446098431612
$ws.Range("O22").Value = 4.944656883147729

$ws.Range("B23").Value = 0.6278170902564
$ws.Range("C23").Value = 0.1317719470093834
$ws.Range("D23").Value = 0.7236935753567764
$ws.Range("E23").Value = 0.2998146715395222
$ws.Range("G23").Value = 0.00248810984962514
$ws.Range("I23").Value = 0.9485541953680681
$ws.Range("J23").Value = 0.1599410469224551
$ws.Range("K23").Value = 0.7336297528378566
$ws.Range("O23").Value = 4.945162946951427

$ws.Range("B24").Value = 0.5587476595956389
$ws.Range("C24").Value = 0.1177033709168143
$ws.Range("D24").Value = 0.7010408589939345
$ws.Range("E24").Value = 0.2889195450554141
$ws.Range("G24").Value = 0.002492589270665905
$ws.Range("I24").Value = 0.9554809481666737
$ws.Range("J24").Value = 0.1528832348669056
$ws.Range("K24").Value = 0.6534967411631101
$ws.Range("O24").Value = 4.951583120545763

$ws.Range("B25").Value = 0.484636245422962
$ws.Range("C25").Value = 0.1026262716961952
$ws.Range("D25").Value = 0.6778836663577295
$ws.Range("E25").Value = 0.2777151184420035
$ws.Range("G25").Value = 0.002497783461363605
$ws.Range("I25").Value = 0.9647447898083072
$ws.Range("J25").Value = 0.1455745388420553
$ws.Range("K25").Value = 0.5675381560945425
$ws.Range("O25").Value = 4.968110786634924
